$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date in A1 (one month later: 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update the three price values in column D
$ws.Range("D14").Value = 1266.597
$ws.Range("D15").Value = 1546.566
$ws.Range("D16").Value = 1817.002
